$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 111.8
$ws.Range("I11").Value = 111.8
$ws.Range("K11").Value = 111.8
$ws.Range("M11").Value = 28.2

$ws.Range("H29").Value = 928.5714
$ws.Range("J29").Value = 992.3077
$ws.Range("L29").Value = 2976.9231
$ws.Range("N29").Value = -3538.9231

$ws.Range("H41").Value = 624.8333
$ws.Range("I41").Value = 599.5
$ws.Range("K41").Value = 599.5
$ws.Range("M41").Value = -159.5

$ws.Range("H43").Value = 998.3333
$ws.Range("J43").Value = 997.5
$ws.Range("L43").Value = 997.5
$ws.Range("N43").Value = -1135.5

$ws.Range("H58").Value = 1494.4445
$ws.Range("I58").Value = 1458.1666
$ws.Range("J58").Value = 1567
$ws.Range("K58").Value = 4374.4998
$ws.Range("L58").Value = 4701
$ws.Range("M58").Value = -4224.4998
$ws.Range("N58").Value = -5001

$ws.Range("H76").Value = 1774.75
$ws.Range("I76").Value = 1366.3334
$ws.Range("K76").Value = 1366.3334
$ws.Range("M76").Value = -1051.3334

$ws.Range("H79").Value = 1774.75
$ws.Range("I79").Value = 1366.3334
$ws.Range("K79").Value = 1366.3334
$ws.Range("M79").Value = -274.3334

$ws.Range("H92").Value = 790
$ws.Range("I92").Value = 856.7857
$ws.Range("K92").Value = 856.7857
$ws.Range("M92").Value = 391.2143

$ws.Range("H103").Value = 533.3333
$ws.Range("I103").Value = 350
$ws.Range("J103").Value = 900
$ws.Range("K103").Value = 1050
$ws.Range("L103").Value = 2700
$ws.Range("M103").Value = -464
$ws.Range("N103").Value = -3872

$ws.Range("H118").Value = 461.33334
$ws.Range("I118").Value = 393.6
$ws.Range("K118").Value = 1180.8
$ws.Range("M118").Value = 476.1999999999998

$ws.Range("H131").Value = 2455.7144
$ws.Range("I131").Value = 1047.5
$ws.Range("K131").Value = 3142.5
$ws.Range("M131").Value = 1897.5

$ws.Range("H135").Value = 1182.6666
$ws.Range("I135").Value = 631.2727
$ws.Range("J135").Value = 2699
$ws.Range("K135").Value = 5681.454299999999
$ws.Range("L135").Value = 24291
$ws.Range("M135").Value = -3146.454299999999
$ws.Range("N135").Value = -29361

$ws.Range("H137").Value = 2936.3157
$ws.Range("I137").Value = 2664.1177
$ws.Range("K137").Value = 7992.353099999999
$ws.Range("M137").Value = -5442.353099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6970.8066
$ws.Range("I32").Value = 6970.8066
$ws.Range("K32").Value = 6970.8066
$ws.Range("M32").Value = -6683.8066

$ws.Range("H74").Value = 27876.125
$ws.Range("I74").Value = 26144.285
$ws.Range("K74").Value = 26144.285
$ws.Range("M74").Value = -25270.285

$ws.Range("H77").Value = 27876.125
$ws.Range("I77").Value = 26144.285
$ws.Range("K77").Value = 130721.425
$ws.Range("M77").Value = -126353.425

$ws.Range("H102").Value = 3010.4546
$ws.Range("I102").Value = 2811.5
$ws.Range("K102").Value = 2811.5
$ws.Range("M102").Value = -1189.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3657.6667
$ws.Range("I94").Value = 3081.25
$ws.Range("K94").Value = 3081.25
$ws.Range("M94").Value = -2630.25

$ws.Range("H105").Value = 2844.8462
$ws.Range("I105").Value = 2844.8462
$ws.Range("K105").Value = 2844.8462
$ws.Range("M105").Value = -1097.8462

$ws.Range("H134").Value = 3049.6667
$ws.Range("I134").Value = 2859.6
$ws.Range("K134").Value = 8578.799999999999
$ws.Range("M134").Value = -6043.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3599.2
$ws.Range("I132").Value = 2666
$ws.Range("K132").Value = 7998
$ws.Range("M132").Value = -5468

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4199.5
$ws.Range("J75").Value = 3239.4
$ws.Range("L75").Value = 9718.200000000001
$ws.Range("N75").Value = -11714.2

$ws.Range("H78").Value = 4199.5
$ws.Range("J78").Value = 3239.4
$ws.Range("L78").Value = 29154.6
$ws.Range("N78").Value = -39138.60000000001

$ws.Range("H87").Value = 6839.6665
$ws.Range("I87").Value = 5756.5
$ws.Range("K87").Value = 17269.5
$ws.Range("M87").Value = -16021.5

$ws.Range("H90").Value = 6839.6665
$ws.Range("I90").Value = 5756.5
$ws.Range("K90").Value = 51808.5
$ws.Range("M90").Value = -45568.5

$ws.Range("H107").Value = 751.5
$ws.Range("I107").Value = 503
$ws.Range("J107").Value = 801.2
$ws.Range("K107").Value = 1509
$ws.Range("L107").Value = 2403.6
$ws.Range("M107").Value = 411
$ws.Range("N107").Value = -6243.6

$ws.Range("H114").Value = 1916.5
$ws.Range("I114").Value = 1749.5
$ws.Range("K114").Value = 5248.5
$ws.Range("M114").Value = -1994.5

$ws.Range("H122").Value = 3168.9546
$ws.Range("I122").Value = 1375
$ws.Range("J122").Value = 3348.35
$ws.Range("K122").Value = 12375
$ws.Range("L122").Value = 30135.15
$ws.Range("M122").Value = -9925
$ws.Range("N122").Value = -35035.14999999999

$ws.Range("H128").Value = 199998
$ws.Range("I128").Value = 199998
$ws.Range("K128").Value = 599994
$ws.Range("M128").Value = -595014

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2299.375
$ws.Range("I7").Value = 2342.1428
$ws.Range("K7").Value = 2342.1428
$ws.Range("M7").Value = -2230.1428

$ws.Range("H16").Value = 259.5
$ws.Range("I16").Value = 299
$ws.Range("J16").Value = 220
$ws.Range("K16").Value = 299
$ws.Range("L16").Value = 220
$ws.Range("M16").Value = -129
$ws.Range("N16").Value = -560

$ws.Range("H61").Value = 15199.8
$ws.Range("I61").Value = 13999.75
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 13999.75
$ws.Range("L61").Value = 20000
$ws.Range("M61").Value = -13797.75
$ws.Range("N61").Value = -20404

$ws.Range("I93").Value = 3451
$ws.Range("J93").Value = 3467.3333
$ws.Range("K93").Value = 3451
$ws.Range("L93").Value = 3467.3333
$ws.Range("M93").Value = -2203
$ws.Range("N93").Value = -5963.3333

$ws.Range("H113").Value = 15199.8
$ws.Range("I113").Value = 13999.75
$ws.Range("J113").Value = 20000
$ws.Range("K113").Value = 13999.75
$ws.Range("L113").Value = 20000
$ws.Range("M113").Value = -11829.75
$ws.Range("N113").Value = -24340

$ws.Range("H122").Value = 3393.7273
$ws.Range("I122").Value = 3625.6667
$ws.Range("K122").Value = 10877.0001
$ws.Range("M122").Value = -8427.000100000001

$ws.Range("H126").Value = 2299.375
$ws.Range("I126").Value = 2342.1428
$ws.Range("K126").Value = 7026.428400000001
$ws.Range("M126").Value = -4556.428400000001

$ws.Range("H136").Value = 4666.3335
$ws.Range("I136").Value = 4666.3335
$ws.Range("K136").Value = 13999.0005
$ws.Range("M136").Value = -11449.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1322.8
$ws.Range("I132").Value = 737
$ws.Range("K132").Value = 2211
$ws.Range("M132").Value = 319

$ws.Range("H136").Value = 2159.25
$ws.Range("I136").Value = 1958.56
$ws.Range("K136").Value = 5875.68
$ws.Range("M136").Value = -3325.68
